$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr14 = New-Object "object[,]" 1,27
$arr14[0,0] = 6772175
$arr14[0,1] = 'Colombia Primera A'
$arr14[0,2] = 45094.91666666666
$arr14[0,3] = 'Atletico Nacional Medellin'
$arr14[0,4] = 'Deportivo Pasto'
$arr14[0,5] = 3
$arr14[0,6] = 2
$arr14[0,7] = 'H'
$arr14[0,8] = 1.666
$arr14[0,9] = 3.75
$arr14[0,10] = 4.5
$arr14[0,11] = 1.8
$arr14[0,12] = 3.6
$arr14[0,13] = 5
$arr14[0,14] = -0.75
$arr14[0,15] = 2
$arr14[0,16] = 1.85
$arr14[0,17] = 2.25
$arr14[0,18] = 1.85
$arr14[0,19] = 2
$arr14[0,20] = 0.8
$arr14[0,21] = -1
$arr14[0,22] = -1
$arr14[0,23] = 0.5
$arr14[0,24] = -0.5
$arr14[0,25] = 0.8500000000000001
$arr14[0,26] = -1
$ws.Range("B14:AB14").Value = $arr14

$arr15 = New-Object "object[,]" 1,27
$arr15[0,0] = 6772177
$arr15[0,1] = 'Colombia Primera A'
$arr15[0,2] = 45094.91666666666
$arr15[0,3] = 'Aguilas Doradas'
$arr15[0,4] = 'Alianza Petrolera'
$arr15[0,5] = 1
$arr15[0,6] = 1
$arr15[0,7] = 'D'
$arr15[0,8] = 2.15
$arr15[0,9] = 3.3
$arr15[0,10] = 3.5
$arr15[0,11] = 2.2
$arr15[0,12] = 3.5
$arr15[0,13] = 3.2
$arr15[0,14] = -0.25
$arr15[0,15] = 1.9
$arr15[0,16] = 1.9
$arr15[0,17] = 2.75
$arr15[0,18] = 1.95
$arr15[0,19] = 1.85
$arr15[0,20] = -1
$arr15[0,21] = 2.5
$arr15[0,22] = -1
$arr15[0,23] = -0.5
$arr15[0,24] = 0.45
$arr15[0,25] = -1
$arr15[0,26] = 0.8500000000000001
$ws.Range("B15:AB15").Value = $arr15

$arr208 = New-Object "object[,]" 1,27
$arr208[0,0] = 7404212
$arr208[0,1] = 'Colombia Primera A'
$arr208[0,2] = 45237.89583333334
$arr208[0,3] = 'Envigado FC'
$arr208[0,4] = 'Deportivo Pasto'
$arr208[0,5] = 1
$arr208[0,6] = 1
$arr208[0,7] = 'D'
$arr208[0,8] = 2.6
$arr208[0,9] = 2.875
$arr208[0,10] = 2.8
$arr208[0,11] = 2.8
$arr208[0,12] = 3.2
$arr208[0,13] = 2.625
$arr208[0,14] = 0
$arr208[0,15] = 1.975
$arr208[0,16] = 1.875
$arr208[0,17] = 2.5
$arr208[0,18] = 2.025
$arr208[0,19] = 1.825
$arr208[0,20] = -1
$arr208[0,21] = 2.2
$arr208[0,22] = -1
$arr208[0,23] = 0
$arr208[0,24] = 0
$arr208[0,25] = -1
$arr208[0,26] = 0.825
$ws.Range("B208:AB208").Value = $arr208

$arr212 = New-Object "object[,]" 1,27
$arr212[0,0] = 7404216
$arr212[0,1] = 'Colombia Primera A'
$arr212[0,2] = 45237.89583333334
$arr212[0,3] = 'Independiente Santa Fe'
$arr212[0,4] = 'Once Caldas'
$arr212[0,5] = 0
$arr212[0,6] = 1
$arr212[0,7] = 'A'
$arr212[0,8] = 1.85
$arr212[0,9] = 3.1
$arr212[0,10] = 4.2
$arr212[0,11] = 2.25
$arr212[0,12] = 3.2
$arr212[0,13] = 3.3
$arr212[0,14] = -0.25
$arr212[0,15] = 1.9
$arr212[0,16] = 1.9
$arr212[0,17] = 2.5
$arr212[0,18] = 1.925
$arr212[0,19] = 1.925
$arr212[0,20] = -1
$arr212[0,21] = -1
$arr212[0,22] = 2.3
$arr212[0,23] = -1
$arr212[0,24] = 0.8999999999999999
$arr212[0,25] = -1
$arr212[0,26] = 0.925
$ws.Range("B212:AB212").Value = $arr212

$arr213 = New-Object "object[,]" 1,27
$arr213[0,0] = 7404260
$arr213[0,1] = 'Colombia Primera A'
$arr213[0,2] = 45238.89583333334
$arr213[0,3] = 'Atletico Nacional Medellin'
$arr213[0,4] = 'Deportes Tolima'
$arr213[0,5] = 2
$arr213[0,6] = 3
$arr213[0,7] = 'A'
$arr213[0,8] = 2
$arr213[0,9] = 3.25
$arr213[0,10] = 3.5
$arr213[0,11] = 1.75
$arr213[0,12] = 3.6
$arr213[0,13] = 4.75
$arr213[0,14] = -0.75
$arr213[0,15] = 2
$arr213[0,16] = 1.8
$arr213[0,17] = 2.5
$arr213[0,18] = 2
$arr213[0,19] = 1.8
$arr213[0,20] = -1
$arr213[0,21] = -1
$arr213[0,22] = 3.75
$arr213[0,23] = -1
$arr213[0,24] = 0.8
$arr213[0,25] = 1
$arr213[0,26] = -1
$ws.Range("B213:AB213").Value = $arr213

$arr215 = New-Object "object[,]" 1,27
$arr215[0,0] = 7404219
$arr215[0,1] = 'Colombia Primera A'
$arr215[0,2] = 45238.89583333334
$arr215[0,3] = 'Union Magdalena'
$arr215[0,4] = 'Independiente Medellin'
$arr215[0,5] = 0
$arr215[0,6] = 4
$arr215[0,7] = 'A'
$arr215[0,8] = 3
$arr215[0,9] = 3.1
$arr215[0,10] = 2.3
$arr215[0,11] = 3.6
$arr215[0,12] = 3.4
$arr215[0,13] = 2.1
$arr215[0,14] = 0.25
$arr215[0,15] = 2.025
$arr215[0,16] = 1.775
$arr215[0,17] = 2.5
$arr215[0,18] = 1.85
$arr215[0,19] = 1.95
$arr215[0,20] = -1
$arr215[0,21] = -1
$arr215[0,22] = 1.1
$arr215[0,23] = -1
$arr215[0,24] = 0.7749999999999999
$arr215[0,25] = 0.8500000000000001
$arr215[0,26] = -1
$ws.Range("B215:AB215").Value = $arr215

$arr216 = New-Object "object[,]" 1,27
$arr216[0,0] = 7404522
$arr216[0,1] = 'Colombia Primera A'
$arr216[0,2] = 45238.89583333334
$arr216[0,3] = 'La Equidad'
$arr216[0,4] = 'Millonarios'
$arr216[0,5] = 2
$arr216[0,6] = 1
$arr216[0,7] = 'H'
$arr216[0,8] = 2.4
$arr216[0,9] = 3.1
$arr216[0,10] = 2.875
$arr216[0,11] = 2.1
$arr216[0,12] = 3.1
$arr216[0,13] = 3.8
$arr216[0,14] = -0.25
$arr216[0,15] = 1.75
$arr216[0,16] = 2.05
$arr216[0,17] = 2
$arr216[0,18] = 1.85
$arr216[0,19] = 1.95
$arr216[0,20] = 1.1
$arr216[0,21] = -1
$arr216[0,22] = -1
$arr216[0,23] = 0.75
$arr216[0,24] = -1
$arr216[0,25] = 0.8500000000000001
$arr216[0,26] = -1
$ws.Range("B216:AB216").Value = $arr216

$arr240 = New-Object "object[,]" 1,27
$arr240[0,0] = 7528135
$arr240[0,1] = 'Colombia Primera A'
$arr240[0,2] = 45266.92708333334
$arr240[0,3] = 'Independiente Medellin'
$arr240[0,4] = 'America de Cali'
$arr240[0,5] = 2
$arr240[0,6] = 1
$arr240[0,7] = 'H'
$arr240[0,8] = 2.15
$arr240[0,9] = 3.3
$arr240[0,10] = 3.4
$arr240[0,11] = 2.375
$arr240[0,12] = 3.3
$arr240[0,13] = 3.1
$arr240[0,14] = -0.25
$arr240[0,15] = 2
$arr240[0,16] = 1.8
$arr240[0,17] = 2.5
$arr240[0,18] = 1.975
$arr240[0,19] = 1.825
$arr240[0,20] = 1.375
$arr240[0,21] = -1
$arr240[0,22] = -1
$arr240[0,23] = 1
$arr240[0,24] = -1
$arr240[0,25] = 0.9750000000000001
$arr240[0,26] = -1
$ws.Range("B240:AB240").Value = $arr240

$arr241 = New-Object "object[,]" 1,27
$arr241[0,0] = 7528603
$arr241[0,1] = 'Colombia Primera A'
$arr241[0,2] = 45266.92708333334
$arr241[0,3] = 'Junior'
$arr241[0,4] = 'Deportes Tolima'
$arr241[0,5] = 4
$arr241[0,6] = 2
$arr241[0,7] = 'H'
$arr241[0,8] = 1.95
$arr241[0,9] = 3.3
$arr241[0,10] = 4
$arr241[0,11] = 1.909
$arr241[0,12] = 3.75
$arr241[0,13] = 3.8
$arr241[0,14] = -0.5
$arr241[0,15] = 1.9
$arr241[0,16] = 1.9
$arr241[0,17] = 2.5
$arr241[0,18] = 1.85
$arr241[0,19] = 1.95
$arr241[0,20] = 0.909
$arr241[0,21] = -1
$arr241[0,22] = -1
$arr241[0,23] = 0.8999999999999999
$arr241[0,24] = -1
$arr241[0,25] = 0.8500000000000001
$arr241[0,26] = -1
$ws.Range("B241:AB241").Value = $arr241

$arr424 = New-Object "object[,]" 1,27
$arr424[0,0] = 7658989
$arr424[0,1] = 'Colombia Primera A'
$arr424[0,2] = 45410.79166666666
$arr424[0,3] = 'Jaguares de Cordoba'
$arr424[0,4] = 'Independiente Santa Fe'
$arr424[0,5] = 1
$arr424[0,6] = 0
$arr424[0,7] = 'H'
$arr424[0,8] = 3
$arr424[0,9] = 3.2
$arr424[0,10] = 2.3
$arr424[0,11] = 3.4
$arr424[0,12] = 3.6
$arr424[0,13] = 2.05
$arr424[0,14] = 0.25
$arr424[0,15] = 2
$arr424[0,16] = 1.8
$arr424[0,17] = 2.5
$arr424[0,18] = 1.8
$arr424[0,19] = 2
$arr424[0,20] = 2.4
$arr424[0,21] = -1
$arr424[0,22] = -1
$arr424[0,23] = 1
$arr424[0,24] = -1
$arr424[0,25] = -1
$arr424[0,26] = 1
$ws.Range("B424:AB424").Value = $arr424

$arr425 = New-Object "object[,]" 1,27
$arr425[0,0] = 7658915
$arr425[0,1] = 'Colombia Primera A'
$arr425[0,2] = 45410.79166666666
$arr425[0,3] = 'Once Caldas'
$arr425[0,4] = 'America de Cali'
$arr425[0,5] = 0
$arr425[0,6] = 0
$arr425[0,7] = 'D'
$arr425[0,8] = 2.3
$arr425[0,9] = 3
$arr425[0,10] = 3.1
$arr425[0,11] = 2.3
$arr425[0,12] = 3.2
$arr425[0,13] = 3.3
$arr425[0,14] = -0.25
$arr425[0,15] = 1.975
$arr425[0,16] = 1.875
$arr425[0,17] = 2.25
$arr425[0,18] = 2.025
$arr425[0,19] = 1.825
$arr425[0,20] = -1
$arr425[0,21] = 2.2
$arr425[0,22] = -1
$arr425[0,23] = -0.5
$arr425[0,24] = 0.4375
$arr425[0,25] = -1
$arr425[0,26] = 0.825
$ws.Range("B425:AB425").Value = $arr425

$arr426 = New-Object "object[,]" 1,27
$arr426[0,0] = 7658985
$arr426[0,1] = 'Colombia Primera A'
$arr426[0,2] = 45410.79166666666
$arr426[0,3] = 'Aguilas Doradas'
$arr426[0,4] = 'Fortaleza'
$arr426[0,5] = 1
$arr426[0,6] = 1
$arr426[0,7] = 'D'
$arr426[0,8] = 1.75
$arr426[0,9] = 3.2
$arr426[0,10] = 5
$arr426[0,11] = 2.05
$arr426[0,12] = 3.2
$arr426[0,13] = 4
$arr426[0,14] = -0.5
$arr426[0,15] = 2.025
$arr426[0,16] = 1.775
$arr426[0,17] = 2
$arr426[0,18] = 1.8
$arr426[0,19] = 2
$arr426[0,20] = -1
$arr426[0,21] = 2.2
$arr426[0,22] = -1
$arr426[0,23] = -1
$arr426[0,24] = 0.7749999999999999
$arr426[0,25] = 0
$arr426[0,26] = 0
$ws.Range("B426:AB426").Value = $arr426

$arr427 = New-Object "object[,]" 1,27
$arr427[0,0] = 7658914
$arr427[0,1] = 'Colombia Primera A'
$arr427[0,2] = 45410.79166666666
$arr427[0,3] = 'La Equidad'
$arr427[0,4] = 'Deportivo Pereira'
$arr427[0,5] = 0
$arr427[0,6] = 2
$arr427[0,7] = 'A'
$arr427[0,8] = 2
$arr427[0,9] = 3.1
$arr427[0,10] = 3.75
$arr427[0,11] = 2.25
$arr427[0,12] = 3.2
$arr427[0,13] = 3.3
$arr427[0,14] = -0.25
$arr427[0,15] = 1.925
$arr427[0,16] = 1.875
$arr427[0,17] = 2
$arr427[0,18] = 1.825
$arr427[0,19] = 1.975
$arr427[0,20] = -1
$arr427[0,21] = -1
$arr427[0,22] = 2.3
$arr427[0,23] = -1
$arr427[0,24] = 0.875
$arr427[0,25] = 0
$arr427[0,26] = 0
$ws.Range("B427:AB427").Value = $arr427

$arr428 = New-Object "object[,]" 1,27
$arr428[0,0] = 7658987
$arr428[0,1] = 'Colombia Primera A'
$arr428[0,2] = 45410.79166666666
$arr428[0,3] = 'Deportivo Cali'
$arr428[0,4] = 'Junior'
$arr428[0,5] = 0
$arr428[0,6] = 0
$arr428[0,7] = 'D'
$arr428[0,8] = 2.7
$arr428[0,9] = 3.25
$arr428[0,10] = 2.4
$arr428[0,11] = 3.2
$arr428[0,12] = 3.1
$arr428[0,13] = 2.4
$arr428[0,14] = 0.25
$arr428[0,15] = 1.8
$arr428[0,16] = 2.05
$arr428[0,17] = 2.25
$arr428[0,18] = 1.975
$arr428[0,19] = 1.875
$arr428[0,20] = -1
$arr428[0,21] = 2.1
$arr428[0,22] = -1
$arr428[0,23] = 0.4
$arr428[0,24] = -0.5
$arr428[0,25] = -1
$arr428[0,26] = 0.875
$ws.Range("B428:AB428").Value = $arr428

$arr429 = New-Object "object[,]" 1,27
$arr429[0,0] = 7658988
$arr429[0,1] = 'Colombia Primera A'
$arr429[0,2] = 45410.79166666666
$arr429[0,3] = 'Envigado FC'
$arr429[0,4] = 'Independiente Medellin'
$arr429[0,5] = 0
$arr429[0,6] = 1
$arr429[0,7] = 'A'
$arr429[0,8] = 4.2
$arr429[0,9] = 3.4
$arr429[0,10] = 1.8
$arr429[0,11] = 5.25
$arr429[0,12] = 3.6
$arr429[0,13] = 1.7
$arr429[0,14] = 0.75
$arr429[0,15] = 1.925
$arr429[0,16] = 1.875
$arr429[0,17] = 2.25
$arr429[0,18] = 1.775
$arr429[0,19] = 2.025
$arr429[0,20] = -1
$arr429[0,21] = -1
$arr429[0,22] = 0.7
$arr429[0,23] = -0.5
$arr429[0,24] = 0.4375
$arr429[0,25] = -1
$arr429[0,26] = 1.025
$ws.Range("B429:AB429").Value = $arr429

$arr430 = New-Object "object[,]" 1,27
$arr430[0,0] = 7736841
$arr430[0,1] = 'Colombia Primera A'
$arr430[0,2] = 45410.79166666666
$arr430[0,3] = 'Atletico Bucaramanga'
$arr430[0,4] = 'Alianza'
$arr430[0,5] = 1
$arr430[0,6] = 0
$arr430[0,7] = 'H'
$arr430[0,8] = 1.666
$arr430[0,9] = 3.5
$arr430[0,10] = 5
$arr430[0,11] = 1.65
$arr430[0,12] = 3.75
$arr430[0,13] = 5.75
$arr430[0,14] = -0.75
$arr430[0,15] = 1.8
$arr430[0,16] = 2
$arr430[0,17] = 2.25
$arr430[0,18] = 1.9
$arr430[0,19] = 1.9
$arr430[0,20] = 0.6499999999999999
$arr430[0,21] = -1
$arr430[0,22] = -1
$arr430[0,23] = 0.4
$arr430[0,24] = -0.5
$arr430[0,25] = -1
$arr430[0,26] = 0.8999999999999999
$ws.Range("B430:AB430").Value = $arr430

$arr431 = New-Object "object[,]" 1,27
$arr431[0,0] = 7658990
$arr431[0,1] = 'Colombia Primera A'
$arr431[0,2] = 45410.79166666666
$arr431[0,3] = 'Millonarios'
$arr431[0,4] = 'Boyaca Chico'
$arr431[0,5] = 3
$arr431[0,6] = 0
$arr431[0,7] = 'H'
$arr431[0,8] = 1.4
$arr431[0,9] = 4.2
$arr431[0,10] = 7
$arr431[0,11] = 1.4
$arr431[0,12] = 4.5
$arr431[0,13] = 8.5
$arr431[0,14] = -1.25
$arr431[0,15] = 1.95
$arr431[0,16] = 1.9
$arr431[0,17] = 2.5
$arr431[0,18] = 1.975
$arr431[0,19] = 1.875
$arr431[0,20] = 0.3999999999999999
$arr431[0,21] = -1
$arr431[0,22] = -1
$arr431[0,23] = 0.95
$arr431[0,24] = -1
$arr431[0,25] = 0.9750000000000001
$arr431[0,26] = -1
$ws.Range("B431:AB431").Value = $arr431

$ws.Range("Q440").Value = 2.025
$ws.Range("R440").Value = 1.825
$ws.Range("T440").Value = 2
$ws.Range("U440").Value = 1.85

$ws.Range("N441").Value = 2.8
$ws.Range("O441").Value = 2.8
$ws.Range("T441").Value = 2.05
$ws.Range("U441").Value = 1.8

$ws.Range("Q442").Value = 2.025
$ws.Range("R442").Value = 1.825

$ws.Range("M443").Value = 1.7
$ws.Range("O443").Value = 5.5
$ws.Range("Q443").Value = 1.925
$ws.Range("R443").Value = 1.925
